$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instrucciones")
$r = $ws.Range("C19")
Write-Host $r.Value()
